$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in new column H, copying the existing header formatting
# (bold, centered, bordered) from G1 so the new column matches the others.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the Save column values (0/1) for data rows 2-15
$saveValues = @(0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
